$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 503.32257
$ws.Range("J17").Value = 503.32257
$ws.Range("L17").Value = 1509.96771
$ws.Range("N17").Value = -1845.96771
$ws.Range("H33").Value = 502.91666
$ws.Range("I33").Value = 505
$ws.Range("K33").Value = 505
$ws.Range("M33").Value = -276
$ws.Range("H76").Value = 5899.909
$ws.Range("I76").Value = 4999.7144
$ws.Range("K76").Value = 4999.7144
$ws.Range("M76").Value = -4684.7144
$ws.Range("H79").Value = 5899.909
$ws.Range("I79").Value = 4999.7144
$ws.Range("K79").Value = 4999.7144
$ws.Range("M79").Value = -3907.7144
$ws.Range("H92").Value = 546.55
$ws.Range("I92").Value = 417.42105
$ws.Range("J92").Value = 3000
$ws.Range("K92").Value = 417.42105
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = 830.5789500000001
$ws.Range("N92").Value = -5496
$ws.Range("H98").Value = 8734.385
$ws.Range("I98").Value = 5736.091
$ws.Range("K98").Value = 5736.091
$ws.Range("M98").Value = -4238.091
$ws.Range("H122").Value = 8734.385
$ws.Range("I122").Value = 5736.091
$ws.Range("K122").Value = 17208.273
$ws.Range("M122").Value = -14758.273
$ws.Range("H132").Value = 15160846
$ws.Range("I132").Value = 17550880
$ws.Range("J132").Value = 23968.666
$ws.Range("K132").Value = 52652640
$ws.Range("L132").Value = 71905.99800000001
$ws.Range("M132").Value = -52650110
$ws.Range("N132").Value = -76965.99800000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2540.7878
$ws.Range("I32").Value = 2519.5615
$ws.Range("K32").Value = 2519.5615
$ws.Range("M32").Value = -2232.5615
$ws.Range("H74").Value = 977.4722
$ws.Range("I74").Value = 614.08
$ws.Range("J74").Value = 1803.3636
$ws.Range("K74").Value = 614.08
$ws.Range("L74").Value = 1803.3636
$ws.Range("M74").Value = 259.92
$ws.Range("N74").Value = -3551.3636
$ws.Range("H77").Value = 977.4722
$ws.Range("I77").Value = 614.08
$ws.Range("J77").Value = 1803.3636
$ws.Range("K77").Value = 3070.4
$ws.Range("L77").Value = 9016.817999999999
$ws.Range("M77").Value = 1297.6
$ws.Range("N77").Value = -17752.818
$ws.Range("H132").Value = 2549.9375
$ws.Range("I132").Value = 2150.0833
$ws.Range("K132").Value = 6450.249899999999
$ws.Range("M132").Value = -3920.249899999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 22728198
$ws.Range("I94").Value = 25000716
$ws.Range("J94").Value = 3010
$ws.Range("K94").Value = 25000716
$ws.Range("L94").Value = 3010
$ws.Range("M94").Value = -25000265
$ws.Range("N94").Value = -3912
$ws.Range("H99").Value = 26316700
$ws.Range("I99").Value = 35715076
$ws.Range("J99").Value = 1242.2
$ws.Range("K99").Value = 35715076
$ws.Range("L99").Value = 1242.2
$ws.Range("M99").Value = -35713578
$ws.Range("N99").Value = -4238.2
$ws.Range("H122").Value = 41000
$ws.Range("J122").Value = 41000
$ws.Range("L122").Value = 41000
$ws.Range("N122").Value = -50800

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1105.7457
$ws.Range("I31").Value = 1083.1404
$ws.Range("K31").Value = 1083.1404
$ws.Range("M31").Value = -788.1404
$ws.Range("H34").Value = 1105.7457
$ws.Range("I34").Value = 1083.1404
$ws.Range("K34").Value = 1083.1404
$ws.Range("M34").Value = -881.1404
$ws.Range("H99").Value = 2194512.5
$ws.Range("I99").Value = 2925417.2
$ws.Range("J99").Value = 1798
$ws.Range("K99").Value = 2925417.2
$ws.Range("L99").Value = 1798
$ws.Range("M99").Value = -2923919.2
$ws.Range("N99").Value = -4794
$ws.Range("H126").Value = 2194512.5
$ws.Range("I126").Value = 2925417.2
$ws.Range("J126").Value = 1798
$ws.Range("K126").Value = 8776251.600000001
$ws.Range("L126").Value = 5394
$ws.Range("M126").Value = -8773781.600000001
$ws.Range("N126").Value = -10334
$ws.Range("H134").Value = 1741.5333
$ws.Range("I134").Value = 1551.9166
$ws.Range("K134").Value = 4655.7498
$ws.Range("M134").Value = -2120.7498

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1212
$ws.Range("H131").Value = 18184626
$ws.Range("J131").Value = 3209.617
$ws.Range("L131").Value = 9628.851000000001
$ws.Range("N131").Value = -19708.851

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2676.7693
$ws.Range("I80").Value = 1771.1428
$ws.Range("J80").Value = 3733.3333
$ws.Range("K80").Value = 1771.1428
$ws.Range("L80").Value = 3733.3333
$ws.Range("M80").Value = -773.1428000000001
$ws.Range("N80").Value = -5729.3333
$ws.Range("H83").Value = 2676.7693
$ws.Range("I83").Value = 1771.1428
$ws.Range("J83").Value = 3733.3333
$ws.Range("K83").Value = 8855.714
$ws.Range("L83").Value = 18666.6665
$ws.Range("M83").Value = -3863.714
$ws.Range("N83").Value = -28650.6665
$ws.Range("H97").Value = 527.8889
$ws.Range("I97").Value = 527.8889
$ws.Range("K97").Value = 527.8889
$ws.Range("M97").Value = -31.88890000000004
$ws.Range("H122").Value = 2766.6667
$ws.Range("I122").Value = 2920
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8760
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6310
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 2064.9033
$ws.Range("I132").Value = 1426.5834
$ws.Range("K132").Value = 4279.7502
$ws.Range("M132").Value = -1749.7502

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 722.375
$ws.Range("I93").Value = 722.375
$ws.Range("K93").Value = 722.375
$ws.Range("M93").Value = 525.625
$ws.Range("H106").Value = 11666.333
$ws.Range("J106").Value = 11666.333
$ws.Range("L106").Value = 11666.333
$ws.Range("N106").Value = -14190.333
$ws.Range("H122").Value = 25760120
$ws.Range("I122").Value = 56668668
$ws.Range("J122").Value = 2996.6667
$ws.Range("K122").Value = 170006004
$ws.Range("L122").Value = 8990.000100000001
$ws.Range("M122").Value = -170003554
$ws.Range("N122").Value = -13890.0001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 447.75
$ws.Range("I81").Value = 433.66666
$ws.Range("J81").Value = 490
$ws.Range("K81").Value = 867.33332
$ws.Range("L81").Value = 980
$ws.Range("M81").Value = 193.66668
$ws.Range("N81").Value = -3102
$ws.Range("H84").Value = 447.75
$ws.Range("I84").Value = 433.66666
$ws.Range("J84").Value = 490
$ws.Range("K84").Value = 4336.6666
$ws.Range("L84").Value = 4900
$ws.Range("M84").Value = 967.3334000000004
$ws.Range("N84").Value = -15508
$ws.Range("H96").Value = 2463
$ws.Range("J96").Value = 2790
$ws.Range("L96").Value = 2790
$ws.Range("N96").Value = -5536
$ws.Range("H105").Value = 34733.332
$ws.Range("I105").Value = 34700
$ws.Range("J105").Value = 34750
$ws.Range("K105").Value = 34700
$ws.Range("L105").Value = 34750
$ws.Range("M105").Value = -31206
$ws.Range("N105").Value = -41738
$ws.Range("H132").Value = 1490.3793
$ws.Range("I132").Value = 1270.6154
$ws.Range("J132").Value = 3395
$ws.Range("K132").Value = 3811.8462
$ws.Range("L132").Value = 10185
$ws.Range("M132").Value = -1281.8462
$ws.Range("N132").Value = -15245
